# Update cryptocurrency price/volume data (Price column D, Volume(1h) column E)
# Values stay as text (matching the original inlineStr cells) even when they
# look numeric, so NumberFormat is forced to "@" (Text) before assignment for
# any Price cell whose new value would otherwise be auto-converted to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.481.98"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.093.17"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.07"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5213"
$ws.Range("E7").Value = "  -3.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4416"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.81"
$ws.Range("E9").Value = "  +15.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08937"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.152"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("E12").Value = "  -4.40%  "
$ws.Range("D13").Value = "2.089.71"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.685"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.685"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.89"
$ws.Range("E16").Value = "  -2.64%  "
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06606"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.16"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.252"
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("D23").Value = "30.534.68"
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.30"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.310"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("D26").Value = "2.342.40"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.25"
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.558"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.73"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.55"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.657"
$ws.Range("E33").Value = "  +6.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.148"
$ws.Range("E34").Value = "  -2.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.901"
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.07"
$ws.Range("E36").Value = "  +3.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02555"
$ws.Range("E37").Value = "  -2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06814"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.476"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.57"
$ws.Range("E40").Value = "  -3.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2255"
$ws.Range("E41").Value = "  -1.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6868"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.253"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6320"
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.195"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.624"
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.238"
$ws.Range("E49").Value = "  +6.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.242"
$ws.Range("E50").Value = "  -4.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.59"
$ws.Range("E51").Value = "  -2.50%  "
